$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "('Avatar of Hope', ['{6}{W}{W}', 'Creature — Avatar', 'If you have 3 or less life, this spell costs {6} less to cast.', 'Flying', 'Avatar of Hope can block any number of creatures.', '4/9'])"

$ws.Range("A2").Value = $newValue

$ws.Rows("3:8").Delete() | Out-Null
